# Auto-generated edit script: updates cryptos list Price (D) and Volume(1h) (E) columns
# to refresh values as part of scheduled GitHub Actions data update.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.202.59'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -1.27%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.659.83'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.86%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.26%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '217.13'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -1.59%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5178'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -2.09%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.25%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2642'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -1.53%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06280'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -1.74%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.84'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -4.72%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07773'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -0.40%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.483'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.21%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.657.33'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -0.95%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.886.07'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -0.90%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5472'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -1.79%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0₅8139'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -2.51%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '65.00'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -1.20%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '26.199.90'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -1.26%  '
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.27%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.618'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -3.16%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '192.42'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.53%  '
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -2.55%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.013'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -4.79%  '
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.32%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '139.39'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.29%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1223'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -3.95%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.282'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -1.80%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '16.17'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -1.01%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.443'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +1.07%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05930'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.274'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -1.58%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.551'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -1.88%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.285'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -4.22%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.585'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -5.93%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9617'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -4.76%  '
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +0.25%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.770'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -0.40%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5678'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -6.83%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.043'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -0.37%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01593'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -1.51%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8542'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -0.46%  '
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +0.29%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.012.87'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -7.31%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '101.08'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +0.37%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.800.41'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -1.00%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0₈111'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -0.80%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '56.53'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -3.38%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.008'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +0.09%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.061'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -0.63%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05167'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -0.62%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4236'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +0.22%  '
